# ReporteDistribucion.xlsx - "Informe" sheet restructuring
# - Move "Manual" column (old J) to sit right before "Duracion" (old M)
# - Insert a new "Diferencia" column right before "Duracion"
# - Update autofilter / dimension / defined names / selection accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informe")

# 1) Move column J (Manual) to just before old column M (Duracion):
#    cut J, insert the cut cells before column M -> Entrada/Salida shift left (K,L -> J,K),
#    Manual ends up at L.
$ws.Columns.Item(10).Cut()
$ws.Columns.Item(13).Insert()

# 2) Insert a new blank column before the (now shifted) old M position (still column 13)
#    and fill its header with the new "Diferencia" label, matching the style of the
#    neighbouring header cells (same style id as Duracion/Manual header, s=5).
$ws.Columns.Item(13).Insert()
$ws.Range("M11").Value = "Diferencia"
$ws.Range("M11").Style = $ws.Range("N11").Style

# Match row-12 (placeholder data row) style of the new cell to its left neighbour too,
# it is a plain body cell like the rest of that row.
$ws.Range("M12").Style = $ws.Range("L12").Style

# 3) Re-apply the AutoFilter over the new full header range
$ws.AutoFilterMode = $false
$ws.Range("A11:U11").AutoFilter()

# 4) Update the hidden _FilterDatabase defined name to track the new autofilter range
for ($i = 1; $i -le $wb.Names.Count(); $i++) {
    $nm = $wb.Names.Item($i)
    if ($nm.Name() -eq "Informe!_FilterDatabase") {
        $nm.RefersTo = "=Informe!`$A`$11:`$U`$11"
    }
}

# 5) Update all the other defined names whose target cell moved because of the
#    column move / insert, and add the brand-new DIFERENCIA name.
function Set-DefinedName($wb, $name, $refersTo) {
    $found = $false
    for ($i = 1; $i -le $wb.Names.Count(); $i++) {
        $nm = $wb.Names.Item($i)
        if ($nm.Name() -eq $name) {
            $nm.RefersTo = $refersTo
            $found = $true
        }
    }
    if (-not $found) {
        $wb.Names.Add($name, $refersTo)
    }
}

Set-DefinedName $wb "CONFIRMACION" "=Informe!`$R`$11"
Set-DefinedName $wb "DISTANCIA" "=Informe!`$Q`$11"
Set-DefinedName $wb "DURACION" "=Informe!`$N`$11"
Set-DefinedName $wb "ENTRADA" "=Informe!`$J`$11"
Set-DefinedName $wb "HORARIO" "=Informe!`$S`$11"
Set-DefinedName $wb "KM" "=Informe!`$O`$11"
Set-DefinedName $wb "MANUAL" "=Informe!`$L`$11"
Set-DefinedName $wb "READ_INACTIVE" "=Informe!`$U`$11"
Set-DefinedName $wb "SALIDA" "=Informe!`$K`$11"
Set-DefinedName $wb "STATE" "=Informe!`$P`$11"
Set-DefinedName $wb "UNREAD_INACTIVE" "=Informe!`$T`$11"
Set-DefinedName $wb "DIFERENCIA" "=Informe!`$M`$11"

# 6) Selection moved to C7
$ws.Range("C7").Select()
